$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 10 (cd2cbb82-... file) gets its own handoff/handback datetimes
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D10").Value = "2016-03-01 09:18:49"
$wsZhCn.Range("G10").Value = "2016-03-01 09:19:32"

# de-de sheet: row 10 (cd2cbb82-... file) gets its own handoff/handback datetimes
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D10").Value = "2016-03-01 09:18:59"
$wsDeDe.Range("G10").Value = "2016-03-01 09:19:50"
